$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 271, pushing the existing
# rows 271:341 down to 275:345 (formatting of the row above is
# carried into the new rows, same as a native Excel row-insert).
$ws.Rows("271:274").Insert()

# Populate the 4 newly inserted rows with the new data.
$newRows = @(
    @{ Row=271; D=44511; H="Sin especificar"; I="Primera"; J=200; K=10000; L=10000; M=10000; N="`$/caja 60 unidades"; O="Limache"; P=167; Q=60 },
    @{ Row=272; D=44511; H="Sin especificar"; I="Primera"; J=400; K=7000;  L=8000;  M=7500;  N="`$/caja 60 unidades"; O="Región de Arica y Parinacota"; P=125; Q=60 },
    @{ Row=273; D=44511; H="Sin especificar"; I="Primera"; J=100; K=10000; L=10000; M=10000; N="`$/caja 60 unidades"; O="Región de O'Higgins"; P=167; Q=60 },
    @{ Row=274; D=44511; H="Sin especificar"; I="Primera"; J=90;  K=8000;  L=8000;  M=8000;  N="`$/caja 60 unidades"; O="Región del Maule"; P=133; Q=60 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino ensalada"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
